{"js": "// Update the date heading and the 25 division-problem answers in the\n// 5x5 practice table. Each value is addressed by its fixed position\n// (title paragraph / table row+col) rather than by searching for the\n// old text, since a couple of the new values collide with other cells'\n// old values (e.g. \"95\u00f75=19, 0\" is both an old value and a new value),\n// which would make a blind global find/replace ambiguous or re-entrant.\n\nconst body = context.document.body;\n\n// 1) Title paragraph: \"2023-08-29 Tuesday\" -> \"2023-08-30 Wednesday\"\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nparagraphs.items[0].insertText(\"2023-08-30 Wednesday\", \"Replace\");\n\n// 2) Table of division problems (5 data rows x 5 columns, each row\n//    separated by 3 blank rows: data rows are at table-row index\n//    0, 4, 8, 12, 16).\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nconst newValues = [\n  [\"33\u00f78=4, 1\", \"83\u00f75=16, 3\", \"74\u00f74=18, 2\", \"48\u00f76=8, 0\", \"79\u00f72=39, 1\"],\n  [\"53\u00f72=26, 1\", \"37\u00f73=12, 1\", \"70\u00f73=23, 1\", \"67\u00f79=7, 4\", \"29\u00f78=3, 5\"],\n  [\"29\u00f74=7, 1\", \"81\u00f77=11, 4\", \"26\u00f74=6, 2\", \"67\u00f77=9, 4\", \"19\u00f75=3, 4\"],\n  [\"25\u00f72=12, 1\", \"55\u00f72=27, 1\", \"93\u00f75=18, 3\", \"93\u00f72=46, 1\", \"25\u00f78=3, 1\"],\n  [\"89\u00f75=17, 4\", \"19\u00f75=3, 4\", \"95\u00f75=19, 0\", \"61\u00f74=15, 1\", \"11\u00f79=1, 2\"],\n];\nconst dataRowIndices = [0, 4, 8, 12, 16];\n\nfor (let r = 0; r < dataRowIndices.length; r++) {\n  const rowIndex = dataRowIndices[r];\n  for (let c = 0; c < 5; c++) {\n    const cell = table.getCell(rowIndex, c);\n    const cellParagraphs = cell.body.paragraphs;\n    cellParagraphs.load(\"items\");\n    await context.sync();\n    cellParagraphs.items[0].insertText(newValues[r][c], \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date heading and the 25 division-problem answers in the\n# 5x5 practice table. Each value is addressed by its fixed position\n# (paragraph index / table cell row+col) rather than by searching for\n# the old text, since a couple of the new values collide with other\n# cells' old values (e.g. \"95\u00f75=19, 0\" is both an old value and a new\n# value), which would make a blind global find/replace ambiguous or\n# re-entrant.\n\n$d = $word.ActiveDocument\n\n# 1) Title paragraph: \"2023-08-29 Tuesday\" -> \"2023-08-30 Wednesday\"\n$d.Paragraphs.Item(1).Range.Text = \"2023-08-30 Wednesday\"\n\n# 2) Table of division problems (Word COM tables are 1-indexed: rows\n#    1..20, columns 1..5). Only rows 1, 5, 9, 13, 17 hold data; the\n#    rows in between are blank spacer rows.\n$table = $d.Tables.Item(1)\n\n$newValues = @(\n    @(\"33\u00f78=4, 1\", \"83\u00f75=16, 3\", \"74\u00f74=18, 2\", \"48\u00f76=8, 0\", \"79\u00f72=39, 1\"),\n    @(\"53\u00f72=26, 1\", \"37\u00f73=12, 1\", \"70\u00f73=23, 1\", \"67\u00f79=7, 4\", \"29\u00f78=3, 5\"),\n    @(\"29\u00f74=7, 1\", \"81\u00f77=11, 4\", \"26\u00f74=6, 2\", \"67\u00f77=9, 4\", \"19\u00f75=3, 4\"),\n    @(\"25\u00f72=12, 1\", \"55\u00f72=27, 1\", \"93\u00f75=18, 3\", \"93\u00f72=46, 1\", \"25\u00f78=3, 1\"),\n    @(\"89\u00f75=17, 4\", \"19\u00f75=3, 4\", \"95\u00f75=19, 0\", \"61\u00f74=15, 1\", \"11\u00f79=1, 2\")\n)\n$dataRows = @(1, 5, 9, 13, 17)\n\nfor ($r = 0; $r -lt $dataRows.Length; $r++) {\n    $rowIndex = $dataRows[$r]\n    for ($c = 1; $c -le 5; $c++) {\n        $table.Cell($rowIndex, $c).Range.Text = $newValues[$r][$c - 1]\n    }\n}\n"}
